$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old content/formatting entirely so we start fresh
$ws.Cells.Clear() | Out-Null

# --- Set new cell values ---
$ws.Range("A1").Value = "This is Lesson 0. If you don't want a lesson 0, fair enough! Just start in row 2"
$ws.Range("A2").Value = "This is Lesson 1. Lesson 1 should be in row 2 etc."
$ws.Range("B2").Value = "If your lessons are over here, set Column Number to 1 in lesson_indexes.json. If they are in another column, set the appropriate column number"
$ws.Range("A3").Value = "Lesson 2"
$ws.Range("A4").Value = "Lesson 3"
$ws.Range("A5").Value = "Lesson 4"
$ws.Range("A6").Value = "You get the idea"
$ws.Range("A7").Value = "Lorum Impusm"
$ws.Range("A8").Value = "Dolor sit amet"

# --- Apply alignment / wrap styles ---
# Style order matters for cellXfs index allocation: B2's center+wrap style
# needs to be created before A1's wrap-only style.
$ws.Range("B2").HorizontalAlignment = -4108   # xlHAlignCenter
$ws.Range("B2").WrapText = $true

$ws.Range("A1").WrapText = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 46.33333333333333
$ws.Columns.Item(2).ColumnWidth = 43.5

# --- Row heights ---
$ws.Range("A1").RowHeight = 30.75
$ws.Range("A2").RowHeight = 43.5

# --- Selection ---
$ws.Range("B2").Select() | Out-Null
